# "fixed verify email for register"
# The "Farid Abdull" / eqarayev4@std.beu.edu.az registrant rows (11-14, 18)
# were bogus/duplicate data; they are replaced with the correct
# "Elmar Qarayev" / elmarqarayev69@gmail.com rows that had been appended
# further down the sheet (rows 15-18), and the now-duplicated trailing
# rows are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Farid Abdull -> Elmar Qarayev, price 66 -> 27 (status stays Pending)
$ws.Range("A11").Value = "Elmar Qarayev"
$ws.Range("B11").Value = "elmarqarayev69@gmail.com"
$ws.Range("C11").Value = 27

# Row 12: Farid Abdull -> Elmar Qarayev, price 156 -> 9, status Rejected -> Accepted
$ws.Range("A12").Value = "Elmar Qarayev"
$ws.Range("B12").Value = "elmarqarayev69@gmail.com"
$ws.Range("C12").Value = 9
$ws.Range("D12").Value = "Accepted"

# Row 13: Farid Abdull -> Elmar Qarayev, price 66 -> 18 (status stays Accepted)
$ws.Range("A13").Value = "Elmar Qarayev"
$ws.Range("B13").Value = "elmarqarayev69@gmail.com"
$ws.Range("C13").Value = 18

# Row 14: Farid Abdull -> Elmar Qarayev, price 156 -> 165 (status stays Accepted)
$ws.Range("A14").Value = "Elmar Qarayev"
$ws.Range("B14").Value = "elmarqarayev69@gmail.com"
$ws.Range("C14").Value = 165

# Row 15: already Elmar Qarayev/27, just fix status Pending -> Accepted
$ws.Range("D15").Value = "Accepted"

# Rows 16-18 were duplicates of the corrected 11-15 rows; remove them.
$ws.Rows("16:18").Delete()
